$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H82").Value = 690
$ws.Range("I82").Value = 690
$ws.Range("K82").Value = 2070
$ws.Range("M82").Value = -1664
$ws.Range("H85").Value = 690
$ws.Range("I85").Value = 690
$ws.Range("K85").Value = 2070
$ws.Range("M85").Value = -666
$ws.Range("H96").Value = 909.2
$ws.Range("I96").Value = 437.5
$ws.Range("J96").Value = 1223.6666
$ws.Range("K96").Value = 1312.5
$ws.Range("L96").Value = 3670.9998
$ws.Range("M96").Value = 60.5
$ws.Range("N96").Value = -6416.9998
$ws.Range("H100").Value = 1407.3334
$ws.Range("I100").Value = 1407.3334
$ws.Range("K100").Value = 1407.3334
$ws.Range("M100").Value = -866.3334
$ws.Range("H116").Value = 10248
$ws.Range("I116").Value = 3747
$ws.Range("K116").Value = 3747
$ws.Range("M116").Value = -305
$ws.Range("H125").Value = 2250
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 4000
$ws.Range("K125").Value = 4500
$ws.Range("L125").Value = 36000
$ws.Range("M125").Value = -2040
$ws.Range("N125").Value = -40920
$ws.Range("H132").Value = 11444.444
$ws.Range("I132").Value = 11444.444
$ws.Range("K132").Value = 34333.33199999999
$ws.Range("M132").Value = -31803.33199999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 4697.125
$ws.Range("I35").Value = 4715.4
$ws.Range("K35").Value = 4715.4
$ws.Range("M35").Value = -4309.4
$ws.Range("H63").Value = 1623.75
$ws.Range("I63").Value = 997.5
$ws.Range("J63").Value = 2250
$ws.Range("K63").Value = 997.5
$ws.Range("L63").Value = 2250
$ws.Range("M63").Value = -311.5
$ws.Range("N63").Value = -3622
$ws.Range("H66").Value = 1623.75
$ws.Range("I66").Value = 997.5
$ws.Range("J66").Value = 2250
$ws.Range("K66").Value = 4987.5
$ws.Range("L66").Value = 11250
$ws.Range("M66").Value = -1555.5
$ws.Range("N66").Value = -18114
$ws.Range("H74").Value = 2470.889
$ws.Range("I74").Value = 2107.6
$ws.Range("J74").Value = 2925
$ws.Range("K74").Value = 2107.6
$ws.Range("L74").Value = 2925
$ws.Range("M74").Value = -1233.6
$ws.Range("N74").Value = -4673
$ws.Range("H77").Value = 2470.889
$ws.Range("I77").Value = 2107.6
$ws.Range("J77").Value = 2925
$ws.Range("K77").Value = 10538
$ws.Range("L77").Value = 14625
$ws.Range("M77").Value = -6170
$ws.Range("N77").Value = -23361
$ws.Range("H132").Value = 6825.375
$ws.Range("I132").Value = 4995.25
$ws.Range("J132").Value = 8655.5
$ws.Range("K132").Value = 14985.75
$ws.Range("L132").Value = 25966.5
$ws.Range("M132").Value = -12455.75
$ws.Range("N132").Value = -31026.5
$ws.Range("H141").Value = 192500
$ws.Range("J141").Value = 192500
$ws.Range("L141").Value = 192500
$ws.Range("N141").Value = -202860
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 75
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 73
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H86").Value = 398.33334
$ws.Range("I86").Value = 398.33334
$ws.Range("K86").Value = 398.33334
$ws.Range("M86").Value = 724.66666
$ws.Range("H89").Value = 398.33334
$ws.Range("I89").Value = 398.33334
$ws.Range("K89").Value = 1991.6667
$ws.Range("M89").Value = 3624.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 691.44446
$ws.Range("I22").Value = 624.8
$ws.Range("J22").Value = 774.75
$ws.Range("K22").Value = 624.8
$ws.Range("L22").Value = 774.75
$ws.Range("M22").Value = -274.8
$ws.Range("N22").Value = -1474.75
$ws.Range("H134").Value = 1878
$ws.Range("I134").Value = 1878
$ws.Range("K134").Value = 5634
$ws.Range("M134").Value = -3099
$ws.Range("H137").Value = 15000
$ws.Range("I137").Value = 15000
$ws.Range("K137").Value = 15000
$ws.Range("M137").Value = -9900
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1375
$ws.Range("I8").Value = 1375
$ws.Range("K8").Value = 4125
$ws.Range("M8").Value = -3986
$ws.Range("H11").Value = 282.83334
$ws.Range("I11").Value = 99
$ws.Range("J11").Value = 374.75
$ws.Range("K11").Value = 297
$ws.Range("L11").Value = 1124.25
$ws.Range("M11").Value = -157
$ws.Range("N11").Value = -1404.25
$ws.Range("H25").Value = 300.22223
$ws.Range("I25").Value = 80.59999999999999
$ws.Range("J25").Value = 574.75
$ws.Range("K25").Value = 241.8
$ws.Range("L25").Value = 1724.25
$ws.Range("M25").Value = -72.79999999999998
$ws.Range("N25").Value = -2062.25
$ws.Range("H30").Value = 300.22223
$ws.Range("I30").Value = 80.59999999999999
$ws.Range("J30").Value = 574.75
$ws.Range("K30").Value = 241.8
$ws.Range("L30").Value = 1724.25
$ws.Range("M30").Value = -139.8
$ws.Range("N30").Value = -1928.25
$ws.Range("H39").Value = 435.08334
$ws.Range("J39").Value = 452.625
$ws.Range("L39").Value = 1357.875
$ws.Range("N39").Value = -1945.875
$ws.Range("H46").Value = 204
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 204
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 612
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -794
$ws.Range("H50").Value = 105
$ws.Range("I50").Value = 105
$ws.Range("K50").Value = 315
$ws.Range("M50").Value = 166
$ws.Range("H53").Value = 105
$ws.Range("I53").Value = 105
$ws.Range("K53").Value = 315
$ws.Range("M53").Value = 166
$ws.Range("H104").Value = 3726
$ws.Range("I104").Value = 3726
$ws.Range("K104").Value = 11178
$ws.Range("M104").Value = -8557
$ws.Range("H139").Value = 2384.75
$ws.Range("J139").Value = 3000
$ws.Range("L139").Value = 9000
$ws.Range("N139").Value = -19280
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4711
$ws.Range("I126").Value = 4800
$ws.Range("K126").Value = 14400
$ws.Range("M126").Value = -11930
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 800
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1390
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 800
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 800
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1014
$ws.Range("H40").Value = 2329.1667
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2364
$ws.Range("H46").Value = 990
$ws.Range("I46").Value = 990
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 990
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -802
$ws.Range("N46").ClearContents()
$ws.Range("H122").Value = 901.5
$ws.Range("I122").Value = 901.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2704.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -254.5
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 479949
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 52089
$ws.Range("I58").Value = 52084
$ws.Range("K58").Value = 52084
$ws.Range("M58").Value = -51776
$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 50000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 150000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -159360
$ws.Range("H81").Value = 950
$ws.Range("I81").Value = 900
$ws.Range("K81").Value = 1800
$ws.Range("M81").Value = -739
$ws.Range("H84").Value = 950
$ws.Range("I84").Value = 900
$ws.Range("K84").Value = 9000
$ws.Range("M84").Value = -3696
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1553.5714
$ws.Range("I122").Value = 979.1667
$ws.Range("K122").Value = 2937.5001
$ws.Range("M122").Value = -487.5001000000002
$ws.Range("H126").Value = 3714.2856
$ws.Range("I126").Value = 3400
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 10200
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -7730
$ws.Range("N126").Value = -18440
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 1699.8572
$ws.Range("I136").Value = 1699.8572
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5099.571599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2549.571599999999
$ws.Range("N136").ClearContents()
